# GoodInfo_v2 - 2021-12-16 未完成
# Update the daily open-position snapshot figures (current price, stock
# value, unrealized P&L, estimated return, and a couple of fee/tax counts)
# for each holding row on the TodayOpenPosition sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TodayOpenPosition")

# Row -> column letter -> new value
$updates = @{
    2  = @{ F = 60.8;   G = 1757; H = 260;  I = 0.1737 }
    3  = @{ F = 282.5;  G = 1689; H = 95;   I = 0.0596; K = 5 }
    4  = @{ F = 81;     G = 3310; H = 483;  I = 0.1709; K = 9 }
    5  = @{ F = 225;    G = 1570; H = 89;   I = 0.0601 }
    6  = @{ F = 42.75;  G = 1363; H = -129; I = -0.0865 }
    7  = @{ F = 278;    G = 1662; H = 173;  I = 0.1162 }
    8  = @{ F = 302;    G = 3612; H = 691;  I = 0.2366 }
    9  = @{ F = 320;    G = 1595; H = 177;  I = 0.1248 }
    10 = @{ F = 44.15;  G = 1628; H = 147;  I = 0.0993 }
    11 = @{ F = 208;    G = 1866; H = 335;  I = 0.2188 }
    12 = @{ F = 17.55;  G = 1486; H = 6;    I = 0.0041 }
    13 = @{ F = 143;    G = 1853; H = 351;  I = 0.2337 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
